# Update gh-pages to output generated at 456a3b4
# Refresh "want to go" counts / sold-out labels across sheets, and
# update the "全部类型" (All types) sheet's 演出 (performance) listing:
# a new event was added (萤火虫动漫游戏嘉年华) on 2024-07-19 which
# pushes the previously adjacent rows (冰兔2024, 跨越二次元, 昨日重现,
# 燃动!!高梨康治) down by one row, while the old "浪漫古典II" concert
# entry drops out of the list (net row count unchanged).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 77
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = "暂时售罄"
$ws.Range("F10").Value = 8184
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 519
$ws.Range("F13").Value = 678
$ws.Range("F14").Value = 539
$ws.Range("F15").Value = 102
$ws.Range("F18").Value = 630
$ws.Range("F19").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 16
$ws.Range("F31").Value = 368
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F43").Value = 7
$ws.Range("F45").Value = 45
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 188
$ws.Range("F4").Value = 322
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 69
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F16").Value = 158
$ws.Range("F19").Value = 0
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 0
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1718
$ws.Range("F3").Value = 420
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 0
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2024-07-19"
$ws.Range("C11").Value = "广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园"
$ws.Range("D11").Value = "新港东路1000号 保利世贸博览馆"
$ws.Range("E11").Value = "2024.07.19 09:00-07.22 17:00"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "暂时售罄"
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=87210"
$ws.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg"
$ws.Range("C12").Value = "广州·冰兔2024线下live「过去和未来」"
$ws.Range("D12").Value = "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）"
$ws.Range("E12").Value = "2024.07.20 20:00-07.20 22:00"
$ws.Range("G12").Value = 198
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=87546"
$ws.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2024-07-20"
$ws.Range("C13").Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$ws.Range("D13").Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$ws.Range("E13").Value = "2024.07.20 19:30-07.20 21:10"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 280
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"
$ws.Range("C14").Value = "广州·昨日重现——唯美英文经典歌曲演唱会"
$ws.Range("D14").Value = "东风中路299号 广州中山纪念堂"
$ws.Range("E14").Value = "2024.07.21 19:30-07.21 21:30"
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 100
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86802"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202405/DR8AvmXe1716802703006.jpeg"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "2024-07-21"
$ws.Range("C15").Value = "广州·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024"
$ws.Range("D15").Value = "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
$ws.Range("E15").Value = "2024.07.21 14:30-07.21 16:00"
$ws.Range("G15").Value = 280
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=87034"
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202406/LINsP2ui1717741701901.png"
$ws.Range("F16").Value = 8184
$ws.Range("F18").Value = 519
$ws.Range("F21").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 630
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 541
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 304
$ws.Range("F35").Value = 375
$ws.Range("F36").Value = 530
$ws.Range("F37").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 368
$ws.Range("F41").Value = 119
$ws.Range("F42").Value = 808
$ws.Range("F45").Value = 0
$ws.Range("F50").Value = 45
